$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.546.12"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").Value = "2.622.15"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'535.18"
$ws.Range("D6").Value = "'143.13"
$ws.Range("E6").Value = "  +0.97%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.570"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("E9").Value = "  +5.45%  "
$ws.Range("E10").Value = "  -2.18%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("D13").Value = "3.082.93"
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").Value = "58.461.11"
$ws.Range("E14").Value = "  -1.62%  "
$ws.Range("D15").Value = "'20.73"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").Value = "2.619.17"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").Value = "'0.0000132"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("D18").Value = "'4.41"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").Value = "'334.12"
$ws.Range("E19").Value = "  -2.19%  "
$ws.Range("D20").Value = "'10.14"
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").Value = "'6.23"
$ws.Range("E21").Value = "  -1.72%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "'66.31"
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("D24").Value = "'0.415"
$ws.Range("E24").Value = "  +1.53%  "
$ws.Range("E25").Value = "  -1.46%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("D27").Value = "'7.11"
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0735"
$ws.Range("E28").Value = "  -1.62%  "
$ws.Range("B29").Value = "USDe"
$ws.Range("C29").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -2.24%  "
$ws.Range("D31").Value = "'5.87"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").Value = "'18.89"
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("D33").Value = "'151.01"
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("D34").Value = "'3.89"
$ws.Range("E34").Value = "  -2.34%  "
$ws.Range("D35").Value = "'0.849"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'1.10"
$ws.Range("E36").Value = "  -2.21%  "
$ws.Range("D37").Value = "'0.813"
$ws.Range("E37").Value = "  -1.92%  "
$ws.Range("E38").Value = "  -2.90%  "
$ws.Range("D39").Value = "'3.57"
$ws.Range("E39").Value = "  +0.73%  "
$ws.Range("D40").Value = "'281.92"
$ws.Range("E40").Value = "  +3.06%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "'0.594"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("D44").Value = "'18.99"
$ws.Range("E44").Value = "  +2.11%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.0937"
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "'0.0527"
$ws.Range("E46").Value = "  +0.46%  "
$ws.Range("D47").Value = "'0.0224"
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("D48").Value = "1.941.34"
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("D49").Value = "'4.45"
$ws.Range("D50").Value = "'17.92"
$ws.Range("E50").Value = "  -3.57%  "
$ws.Range("D51").Value = "'113.91"
$ws.Range("E51").Value = "  +1.68%  "
